# fix: unique command names in XLSX - prefix protocol name to each step
#
# For every "protocol" worksheet (every sheet except the first 5 overview /
# reference sheets: LanaJourney, NRWaves, PersonalLana, PositiveSpin,
# ReEngagement), prefix each non-empty "Name" cell (column A, rows below the
# header row) with the worksheet's own name followed by a space - unless it
# already begins with that prefix.

$wb = $excel.ActiveWorkbook

# Sheets that must NOT be touched by this change.
$skipNames = @("LanaJourney", "NRWaves", "PersonalLana", "PositiveSpin", "ReEngagement")

foreach ($ws in $wb.Worksheets) {
    if ($skipNames -contains $ws.Name) {
        continue
    }

    $prefix = $ws.Name + " "

    $usedRange = $ws.UsedRange
    $lastRow = $usedRange.Rows.Count

    # Row 1 is the header ("Name", "Text", "Note", "*Guidelines"); data starts
    # on row 2.
    for ($r = 2; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, 1)
        $current = $cell.Value()

        if ($null -eq $current) {
            continue
        }

        $text = [string]$current
        if ($text -eq "") {
            continue
        }

        if ($text.StartsWith($prefix)) {
            continue
        }

        $cell.Value = $prefix + $text
    }
}
